# ----------------------------------------------------------------------
# This script applies the commit's changes to PlayerPerformance_4239.xlsx:
#   1. Insert a new "Player Info" sheet before "ODI Batting" with player
#      metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#   2. In the existing "ODI Batting" and "ODI Bowling" sheets, rename the
#      MATCH_CARD_LINK column to MATCH_CODE and replace the full howstat
#      URL values with just the trailing MatchCode number.
#   3. Append a new "ODI Batting Extra" sheet at the end with additional
#      per-match batting stats (batting position, boundary counts, % of
#      team runs, man-of-the-match flag).
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$xlCenter = -4108
$xlTop = -4160
$xlDown = -4121

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = $xlCenter
    $range.VerticalAlignment = $xlTop
}

function Set-TextCell($cell, $text) {
    # Prefix with an apostrophe so Excel treats numeric-looking values
    # (match codes, percentages, counts, ...) as literal text, matching
    # the inlineStr cells produced by the source scraping script.
    $cell.Value = "'" + $text
}

# ------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the first existing sheet
# ------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Set-HeaderStyle($playerInfo.Range("A1:D1"))

Set-TextCell $playerInfo.Range("A2") "4239"
$playerInfo.Range("B2").Value = "Hashmatullah Shahidi"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ------------------------------------------------------------------
# 2. Update MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" (col D)
#    and "ODI Bowling" (col B), replacing URLs with bare match codes
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingLastRow = $battingSheet.Range("A1").End($xlDown).Row
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Text
    $parts = $url.Split("=")
    $code = $parts[$parts.Length - 1]
    Set-TextCell $cell $code
}

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingLastRow = $bowlingSheet.Range("A1").End($xlDown).Row
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Text
    $parts = $url.Split("=")
    $code = $parts[$parts.Length - 1]
    Set-TextCell $cell $code
}

# ------------------------------------------------------------------
# 3. New "ODI Batting Extra" sheet appended at the end
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"
Set-HeaderStyle($extra.Range("A1:F1"))

$extraData = @(
    @("4315", "4", "9", "0", "34.30%", "NO"),
    @("4323", "4", "1", "0", "6.40%", "NO"),
    @("4326", "", "", "", "", "NO"),
    @("4332", "", "", "", "", "NO"),
    @("4335", "3", "0", "0", "5.50%", "NO"),
    @("4340", "3", "0", "0", "", "NO"),
    @("4446", "4", "8", "1", "31.54%", "NO"),
    @("4448", "4", "0", "0", "3.76%", "NO"),
    @("4525", "", "", "", "", "NO"),
    @("4528", "4", "1", "2", "22.78%", "NO"),
    @("4530", "4", "1", "0", "11.02%", "NO"),
    @("4537", "4", "3", "1", "13.02%", "NO"),
    @("4538", "", "", "", "", "NO"),
    @("4539", "4", "0", "0", "1.04%", "NO"),
    @("4582", "4", "13", "0", "31.88%", "NO"),
    @("4585", "4", "0", "0", "0.44%", "NO"),
    @("4588", "4", "5", "0", "27.74%", "NO"),
    @("4671", "", "", "", "", "NO"),
    @("4674", "4", "2", "0", "12.28%", "NO"),
    @("4675", "", "", "", "", "NO")
)

$row = 2
foreach ($entry in $extraData) {
    Set-TextCell $extra.Cells.Item($row, 1) $entry[0]

    $pos = $entry[1]
    if ($pos -ne "") {
        $extra.Cells.Item($row, 2).Value = [double]$pos
    }

    if ($entry[2] -ne "") {
        Set-TextCell $extra.Cells.Item($row, 3) $entry[2]
    }
    if ($entry[3] -ne "") {
        Set-TextCell $extra.Cells.Item($row, 4) $entry[3]
    }
    if ($entry[4] -ne "") {
        Set-TextCell $extra.Cells.Item($row, 5) $entry[4]
    }

    $extra.Cells.Item($row, 6).Value = $entry[5]

    $row = $row + 1
}

# Keep the first sheet as the active tab, matching the original workbook view
$wb.Worksheets.Item(1).Activate()

